$d = $word.ActiveDocument

# Locate the paragraph that contains the bookmarked heading text
# "Realización de los Casos de Uso de Negocio" and append a formatted
# trailing space run right after it (after the bookmarkEnd), i.e. at
# the very end of that paragraph, before the paragraph mark.
$found = $d.Content.Find.Execute("Realización de los Casos de Uso de Negocio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Realización de los Casos de Uso de Negocio*") {
        $target = $p
    }
}

$r = $target.Range
$r.Collapse(0)
[void]$r.MoveEnd(1, -1)
[void]$r.InsertAfter(" ")
$r.Font.Name = "Tahoma"
$r.Font.NameAscii = "Tahoma"
$r.Font.NameFarEast = "Tahoma"
$r.Font.Italic = $false
$r.Font.Size = 11
$r.LanguageID = 3082
